$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (refreshed data pull) for rows 2-34
$ws.Range("C2").Value = 2779715780.1971998
$ws.Range("C3").Value = 2364740766.4398003
$ws.Range("C4").Value = 2354120786.6111999
$ws.Range("C5").Value = 1647063840.4191999
$ws.Range("C6").Value = 1202672016.6792002
$ws.Range("C7").Value = 1153129978.0010002
$ws.Range("C8").Value = 968327492.69599998
$ws.Range("C9").Value = 648955510.75940001
$ws.Range("C10").Value = 590539773.57160008
$ws.Range("C11").Value = 579670695.75
$ws.Range("C12").Value = 450714236.7712
$ws.Range("C13").Value = 411764926.16729999
$ws.Range("C14").Value = 289204784.35119998
$ws.Range("C15").Value = 252411594.55039999
$ws.Range("C16").Value = 233891214.64499998
$ws.Range("C17").Value = 205750999.59819999
$ws.Range("C18").Value = 138072167.338
$ws.Range("C19").Value = 114976184.516
$ws.Range("C20").Value = 108158853.3408
$ws.Range("C21").Value = 103098592.90099999
$ws.Range("C22").Value = 102258845.814
$ws.Range("C23").Value = 66382476.8719
$ws.Range("C24").Value = 56256615.114200003
$ws.Range("C25").Value = 53401230.173500001
$ws.Range("C26").Value = 52439319.640000001
$ws.Range("C27").Value = 51238979.598899998
$ws.Range("C28").Value = 45036484.104999997
$ws.Range("C29").Value = 33605904.375500001
$ws.Range("C30").Value = 31856096.148400001
$ws.Range("C31").Value = 29982580.704
$ws.Range("C32").Value = 20157133.417599998
$ws.Range("C33").Value = 7540281.0574000003
$ws.Range("C34").Value = 2841504.14

# Update the active selection as recorded in the sheet view
$ws.Range("A2:E34").Select()
